# "Generate Report for Handback"
# Refresh the handback timestamps recorded for the f1426313-...xlf entry
# (last row of each language report) on the zh-cn and de-de sheets.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D8").Value = "2016-02-26 05:02:29"
$wsZhCn.Range("G8").Value = "2016-02-26 05:03:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D8").Value = "2016-02-26 05:02:39"
$wsDeDe.Range("G8").Value = "2016-02-26 05:03:31"
